# Aggiornamento fino a 13/03 (13/05 secondo le date seriali presenti)
# Appende quattro nuove righe di dati (252-255) in coda al foglio,
# replicando lo stesso schema/stile della riga precedente (251).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuovi dati: data seriale Excel, nuovi positivi, somma mobile 7gg, somma mobile 7gg per 100mila abitanti
$newRows = @(
    @(44326, 0, 4, 45.99816007359706),
    @(44327, 0, 3, 34.4986200551978),
    @(44328, 0, 3, 34.4986200551978),
    @(44329, 0, 1, 11.49954001839926)
)

$lastRow = 251

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $targetRow = $lastRow + 1 + $i

    # Copia la formattazione (stile data in colonna A, ecc.) dall'ultima riga esistente
    $srcRange = $ws.Range("A$lastRow" + ":D$lastRow")
    $dstRange = $ws.Range("A$targetRow" + ":D$targetRow")
    $srcRange.Copy($dstRange)

    $vals = $newRows[$i]
    $ws.Cells.Item($targetRow, 1).Value = $vals[0]
    $ws.Cells.Item($targetRow, 2).Value = $vals[1]
    $ws.Cells.Item($targetRow, 3).Value = $vals[2]
    $ws.Cells.Item($targetRow, 4).Value = $vals[3]
}
